$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("B3").Value = 32
$ws.Range("A5").Value = "title"
$ws.Range("B5").Value = 35
$ws.Range("A7").Value = "subtitle"
$ws.Range("B7").Value = 15
$ws.Range("A9").Value = "spacer"
$ws.Range("B9").Value = 20

$ws.Rows.Item(5).RowHeight = 28
$ws.Rows.Item(9).RowHeight = 27

# Title style for A5
$titleRange = $ws.Range("A5")
$titleRange.Font.Bold = $true
$titleRange.Font.Size = 14
$titleRange.Interior.Pattern = -4124
$titleRange.Interior.ThemeColor = 1
$titleRange.Interior.TintAndShade = -0.049989318521683403
$titleRange.Borders.Item(9).LineStyle = 1
$titleRange.Borders.Item(9).Weight = 2
$titleRange.HorizontalAlignment = -4131

# Subtitle style for A7
$subtitleRange = $ws.Range("A7")
$subtitleRange.Font.Italic = $true
$subtitleRange.Font.Size = 8
$subtitleRange.VerticalAlignment = -4108

# Spacer style for A9
$spacerRange = $ws.Range("A9")
$spacerRange.VerticalAlignment = -4108

[void]$ws.Range("C19").Select()

